# 1st changes of mifos to finflux
#
# Insert a new blank column before column N on the "Repayment Schedule"
# sheet, shifting the old N/O/P columns (Late / blank / Outstanding) one
# column to the right (-> O/P/Q), and update the active sheet / selection
# to reflect what was left active after the edit.

$wb = $excel.ActiveWorkbook

$wsSchedule = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new column at N; existing N,O,P data (and the trailing header
# labels "Late" / "Outstanding") shift right to O,P,Q.
$wsSchedule.Columns("N").Insert()

# The "Repayment Schedule" sheet becomes the active sheet/tab, with a new
# selection, while "Transactions" (previously active) automatically loses
# its tabSelected flag since only one sheet can be active at a time.
$wsSchedule.Activate()
$wsSchedule.Range("S6").Select()
